$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF headers, copying the formatting (style) of the
# existing H1 header cell (bold / bordered / centered) onto the new cells.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Data rows 2-13: I column = 1 (constant), J column = same value as H column
for ($r = 2; $r -le 13; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
